$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# Swap values of columns A, B, E, F, G, H between row 2 and row 3
$cols = @("A", "B", "E", "F", "G", "H")

foreach ($col in $cols) {
    $cell2 = $ws.Range($col + "2")
    $cell3 = $ws.Range($col + "3")
    $val2 = $cell2.Value()
    $val3 = $cell3.Value()
    $cell2.Value = $val3
    $cell3.Value = $val2
}
